$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item("testdata_Mean")
$ws2 = $wb.Worksheets.Item("testdata_Mean_results")

# ---------------------------------------------------------------------------
# Sheet "testdata_Mean_results": rebuild with the new confidence-interval
# layout (area, value_sum, value_count, stdev, mean, lowercl, uppercl,
# confidence, statistic, method) for both the 95% and 99.8% methods.
# ---------------------------------------------------------------------------

# Clear old cell values only (keeps the header row's bold style s=3 and the
# sheet's column widths), so stale shared strings for the dropped columns
# (lower95cl/upper95cl/lower998cl/upper998cl/numrecs/total) get pruned on
# save once nothing references them any more.
$ws2.Cells.ClearContents()

# --- Header row: write everything except B1/C1 first, left to right, so new
# shared strings are appended in the exact order the target file expects.
$ws2.Range("A1").Value = "area"
$ws2.Range("D1").Value = "stdev"
$ws2.Range("E1").Value = "mean"
$ws2.Range("F1").Value = "lowercl"
$ws2.Range("G1").Value = "uppercl"
$ws2.Range("H1").NumberFormat = "@"
$ws2.Range("H1").Font.Bold = $true
$ws2.Range("H1").Value = "confidence"
$ws2.Range("I1").Value = "statistic"
$ws2.Range("J1").Value = "method"

# --- Confidence column (H) for every data row next, so "95%" then "99.8%"
# become the next two new shared strings, in that order. Format the whole
# column as text so the percent-like labels aren't coerced into numbers.
$ws2.Columns.Item(8).NumberFormat = "@"
$ws2.Range("H2:H4").Value = "95%"
$ws2.Range("H5:H7").Value = "99.8%"

# --- Statistic column (I): reuses the existing "mean" shared string.
$ws2.Range("I2:I7").Value = "mean"

# --- Method column (J): introduces the last new shared string.
$ws2.Range("J2:J7").Value = "Student's t-distribution"

# --- Area labels (column A) reuse existing shared strings.
$ws2.Range("A2").Value = "Area1"
$ws2.Range("A3").Value = "Area2"
$ws2.Range("A4").Value = "No grouping"
$ws2.Range("A5").Value = "Area1"
$ws2.Range("A6").Value = "Area2"
$ws2.Range("A7").Value = "No grouping"

# --- Numeric data, 95% block (rows 2-4).
$ws2.Range("B2").Value = 429.63476000000009
$ws2.Range("C2").Value = 8
$ws2.Range("D2").Value = 20.975575257709309
$ws2.Range("E2").Value = 53.704345000000011
$ws2.Range("F2").Value = 36.168325241336504
$ws2.Range("G2").Value = 71.240364758663517

$ws2.Range("B3").Value = 102221.33323999999
$ws2.Range("C3").Value = 18
$ws2.Range("D3").Value = 2117.8317161590671
$ws2.Range("E3").Value = 5678.9629577777778
$ws2.Range("F3").Value = 4625.7900224529722
$ws2.Range("G3").Value = 6732.1358931025834

$ws2.Range("B4").Formula = "=SUM(B2:B3)"
$ws2.Range("C4").Formula = "=SUM(C2:C3)"
$ws2.Range("D4").Value = 3171.8018122194453
$ws2.Range("E4").Value = 3948.1141538461534
$ws2.Range("F4").Value = 2666.9956767458489
$ws2.Range("G4").Value = 5229.2326309464579

# --- Numeric data, 99.8% block (rows 5-7).
$ws2.Range("B5").Value = 429.63476000000009
$ws2.Range("C5").Value = 8
$ws2.Range("D5").Value = 20.975575257709309
$ws2.Range("E5").Value = 53.704345000000011
$ws2.Range("F5").Value = 18.216705294788838
$ws2.Range("G5").Value = 89.191984705211183

$ws2.Range("B6").Value = 102221.33323999999
$ws2.Range("C6").Value = 18
$ws2.Range("D6").Value = 2117.8317161590671
$ws2.Range("E6").Value = 5678.9629577777778
$ws2.Range("F6").Value = 3859.0770997295967
$ws2.Range("G6").Value = 7498.8488158259588

$ws2.Range("B7").Formula = "=SUM(B5:B6)"
$ws2.Range("C7").Formula = "=SUM(C5:C6)"
$ws2.Range("D7").Value = 3171.8018122194453
$ws2.Range("E7").Value = 3948.1141538461534
$ws2.Range("F7").Value = 1801.9535385474737
$ws2.Range("G7").Value = 6094.274769144833

# --- Finally the two new header cells, appended last among new strings.
$ws2.Range("B1").Value = "value_sum"
$ws2.Range("C1").Value = "value_count"

# --- Column widths: col 3 (value_count) widens to match col 1/2, col 10
# (method) widens to fit "Student's t-distribution" style labels.
$ws2.Columns.Item(3).ColumnWidth = 11
$ws2.Columns.Item(10).ColumnWidth = 21.166666666666668

# ---------------------------------------------------------------------------
# Sheet "testdata_Mean": header cells now read area / values (they keep the
# same text, just now pointing at the shared strings that survived pruning).
# ---------------------------------------------------------------------------
$ws1.Range("A1").Value = "area"
$ws1.Range("B1").Value = "values"

# --- Selections: sheet1 first, sheet2 last so sheet2 stays the active tab
# (matches tabSelected="1" staying on testdata_Mean_results).
$ws1.Range("C38").Select()
$ws2.Range("F12").Select()
